# Surat Pengunduran Diri template fix
# - Every "${Foo_bar}" style placeholder becomes "${text_Foo_bar}" (the
#   "text_" type-prefix that the templating engine expects), except the
#   "${datetime_...}" one, which already carries its type prefix.
# - "${Nama Anda}" becomes "${text_Nama Anda_nulldesc}" (adds a null
#   description suffix), and the very last occurrence - the signature
#   line at the bottom of the letter - loses its closing "}" while
#   picking up the (hidden) "_GoBack" bookmark that used to sit in the
#   big "Melalui kedatangan surat ini, ..." paragraph.

$d = $word.ActiveDocument

function Replace-All($find, $replace) {
    $r = $d.Content
    $r.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1) "${Nama Anda}" -> "${text_Nama Anda_nulldesc}" (both occurrences;
#    the trailing "}" is stripped back off the very last one afterwards)
Replace-All '${Nama Anda}' '${text_Nama Anda_nulldesc}'

# 2) "${Nama Perusahaan_Perusahaan tempat Anda bekerja}" -> add "text_"
Replace-All '${Nama Perusahaan_Perusahaan tempat Anda bekerja}' '${text_Nama Perusahaan_Perusahaan tempat Anda bekerja}'

# 3) "${Alamat Perusahaan_Tempat perusahaan Anda berada}" -> add "text_"
Replace-All '${Alamat Perusahaan_Tempat perusahaan Anda berada}' '${text_Alamat Perusahaan_Tempat perusahaan Anda berada}'

# 4) "${Kota penulisan surat_Kota tempat surat ditulis}" -> add "text_"
Replace-All '${Kota penulisan surat_Kota tempat surat ditulis}' '${text_Kota penulisan surat_Kota tempat surat ditulis}'

# 5) "${Posisi Anda" (continues "_Jabatan Anda di perusahaan}" in the next
#    run) -> add "text_"
Replace-All '${Posisi Anda' '${text_Posisi Anda'

# 6) The closing "}" on the final signature line ("${text_Nama
#    Anda_nulldesc}" on its own right-aligned paragraph) is dropped, and
#    the "_GoBack" bookmark - previously wrapped around the space right
#    after the big paragraph's date sentence - is re-anchored to the end
#    of that same final paragraph.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$lastText = $lastRange.Text
$trimmed = $lastText.Substring(0, $lastText.Length - 1)
$lastRange.Text = $trimmed

$bm = $d.Bookmarks.Item("_GoBack")
$bmRange = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
